$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999825041397,
    0.99759413077346359,
    0.9921984474598311,
    0.99409703316125653,
    0.98281375739562082,
    0.95538787411995629,
    0.95331861615370395,
    0.94964367894257018,
    0.94748950114919062,
    0.94627208852747313,
    0.9461572579059736,
    0.94619522973527004,
    0.94881209472102501,
    0.95099144357740517,
    0.95423867695041964,
    0.95859699588129366,
    0.95488933015564159,
    0.95378040758919891,
    0.99482482903293812,
    0.98770757176798751,
    0.9863090448630778,
    0.98504453267020842,
    0.98516409503145708,
    0.97214368811185126,
    0.96568675825796657,
    0.95173360556228015,
    0.94688868998081976,
    0.92542119261253153,
    0.91015155038792694,
    0.90358141313051266,
    0.89592759289822332,
    0.89424820269600092,
    0.89372816919394205
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
